$d = $word.ActiveDocument

# Helper: re-merge a run sequence that currently reads as $search (possibly split
# across several runs, e.g. by spell-check proofing marks) into a single run
# reading $replacement. A direct Find/Replace that writes back the *same*
# characters is treated as a no-op by the engine (the runs stay split), so when
# $replacement is textually unchanged we first write a temporary sentinel
# (forcing a genuine content mutation, which collapses the spanned runs into
# one), then rewrite that sentinel down to the final text.
function Merge-Run {
    param(
        $Doc,
        [string]$Search,
        [string]$Replacement,
        [bool]$Hyperlink
    )

    $marker = "@@MRK@@"

    $rng = $Doc.Content
    $found = $rng.Find.Execute($Search, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $Search"
        return
    }

    $rng.Text = $Replacement + $marker

    $rng2 = $Doc.Content
    $found2 = $rng2.Find.Execute($Replacement + $marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found2) {
        Write-Output "MARKER NOT FOUND for: $Search"
        return
    }
    $rng2.Text = $Replacement

    if ($Hyperlink) {
        $rng3 = $Doc.Content
        $found3 = $rng3.Find.Execute($Replacement, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found3) {
            $rng3.Style = "Hyperlink"
        }
    }
}

# 1. Title: prepend ". " before "Evaluating Capabilities..." (new leading run)
$d.Content.Find.Execute("Evaluating Capabilities of current online news scraping software", $true, $false, $false, $false, $false, $true, 1, $false, ". Evaluating Capabilities of current online news scraping software", 2) | Out-Null

# 2. "Two of the most effective ... softwares available are" - merge spell-check-split runs
Merge-Run $d "Two of the most effective open-source news scraping softwares available are" "Two of the most effective open-source news scraping softwares available are"

# 3. "Newspaper3k by Lucas Ou-Yang" (inside hyperlink) - merge spell-check-split runs
Merge-Run $d "Newspaper3k by Lucas Ou-Yang" "Newspaper3k by Lucas Ou-Yang" $true

# 4. "ews-please by Felix Hamborg" (inside hyperlink) - merge spell-check-split runs
Merge-Run $d "ews-please by Felix Hamborg" "ews-please by Felix Hamborg" $true

# 5. "News articles with multiple pages (e.g. chinadaily)" - merge spell-check-split runs
Merge-Run $d "News articles with multiple pages (e.g. chinadaily)" "News articles with multiple pages (e.g. chinadaily)" $false

# 6. "News-please has more informative errors when running on jupyter notebook" - merge spell-check-split runs
Merge-Run $d "News-please has more informative errors when running on jupyter notebook" "News-please has more informative errors when running on jupyter notebook" $false

# 7/8. "Al-Akhbar" (two occurrences in tables) - merge spell-check-split runs
Merge-Run $d "Al-Akhbar" "Al-Akhbar" $false
Merge-Run $d "Al-Akhbar" "Al-Akhbar" $false

# 9. "Kyunghyang Shinmun" - merge spell-check-split runs
Merge-Run $d "Kyunghyang Shinmun" "Kyunghyang Shinmun" $false

# 10. "-Chinadaily article with multiple pages" - merge spell-check-split runs
Merge-Run $d "-Chinadaily article with multiple pages" "-Chinadaily article with multiple pages" $false

# 11. "-The scraper does not correctly parse the chinadaily article with pages..." - merge spell-check-split runs
Merge-Run $d "-The scraper does not correctly parse the chinadaily article with pages (returns blank text even with language specified)" "-The scraper does not correctly parse the chinadaily article with pages (returns blank text even with language specified)" $false

# 12. ", 3 times per article (using timeit)" - merge spell-check-split runs
Merge-Run $d ", 3 times per article (using timeit)" ", 3 times per article (using timeit)" $false

# 13. "Average Time" -> "Average Time (ms)" (new trailing run)
$d.Content.Find.Execute("Average Time", $true, $false, $false, $false, $false, $true, 1, $false, "Average Time (ms)", 2) | Out-Null
